$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Mint the "No Spacing" paragraph style (w:styleId="NoSpacing") used
#    throughout the rewritten email, using the trailing empty paragraph as a
#    scratch anchor so the real content paragraphs are left untouched. Then
#    tune the style definition (BaseStyle/Priority/ParagraphFormat) so the
#    emitted <w:style> matches Word's real "No Spacing" quick style, and put
#    that anchor paragraph back to the Normal style afterwards.
# ---------------------------------------------------------------------------
$anchor = $d.Paragraphs($d.Paragraphs.Count)
$anchor.Style = "No Spacing"
$noSpacing = $anchor.Style
$noSpacing.BaseStyle = ""
$noSpacing.Priority = 1
$noSpacingFmt = $noSpacing.ParagraphFormat
$noSpacingFmt.SpaceAfter = 0
$noSpacingFmt.LineSpacingRule = 0
$anchor2 = $d.Paragraphs($d.Paragraphs.Count)
$anchor2.Style = "Normal"

# ---------------------------------------------------------------------------
# 2. Remove all existing body paragraphs except the trailing empty one
#    (which stays, unstyled, as the document's final paragraph mark).
# ---------------------------------------------------------------------------
$deleteEnd = $d.Paragraphs($d.Paragraphs.Count - 1).Range.End
$d.Range(0, $deleteEnd).Delete()

# ---------------------------------------------------------------------------
# 3. Rebuild the email body as WordprocessingML and insert it in front of
#    the remaining trailing paragraph.
# ---------------------------------------------------------------------------
$bodyXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships"><w:body><w:p><w:pPr><w:pStyle w:val="NoSpacing"/></w:pPr><w:r><w:t>Hi Colby,</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="NoSpacing"/></w:pPr></w:p><w:p><w:pPr><w:pStyle w:val="NoSpacing"/><w:rPr><w:sz w:val="17"/><w:szCs w:val="17"/></w:rPr></w:pPr><w:r><w:t xml:space="preserve">My name is Eduardo Gutarra, and I am reaching out to you because you might happen to know my Master’s thesis supervisor, Prof. Daniel Lemire. I have worked with </w:t></w:r><w:r><w:t xml:space="preserve">EWAH </w:t></w:r><w:r><w:t>a bitmap index</w:t></w:r><w:r><w:t>ing compression technique for</w:t></w:r><w:r><w:t xml:space="preserve"> a job with ION Geophysical and for my master’s degree thesis. I believe you have used EWAH too. I also coauthored in “Reordering Rows for Better Compression” (</w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t xml:space="preserve">see  </w:t></w:r><w:proofErr w:type="gramEnd"/><w:hyperlink r:id="rId4" w:tgtFrame="_blank" w:history="1"><w:r><w:rPr><w:sz w:val="17"/><w:u w:val="single"/></w:rPr><w:t>http://arxiv.org/abs/1207.2189</w:t></w:r></w:hyperlink><w:r><w:rPr><w:sz w:val="17"/><w:szCs w:val="17"/></w:rPr><w:t>)</w:t></w:r><w:r><w:t xml:space="preserve"> for my master’s degree</w:t></w:r><w:r><w:rPr><w:sz w:val="17"/></w:rPr><w:t xml:space="preserve">  </w:t></w:r><w:r><w:t>thesis. It was published in ACM TODS.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="NoSpacing"/><w:rPr><w:sz w:val="17"/><w:szCs w:val="17"/></w:rPr></w:pPr><w:r><w:t>I am currently in the interview process with google and was asked if I know someone that works for them, and I really don’t know anybody personally in google. So, I immediately thought of the possibility of getting to know somebody that is working for google and also has done similar work to mine. Thus, I asked Daniel</w:t></w:r><w:r><w:t xml:space="preserve"> Lemire</w:t></w:r><w:r><w:t xml:space="preserve"> if he knew someone, and he referred me to you.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="NoSpacing"/></w:pPr><w:r><w:t xml:space="preserve">I wanted to ask you for advice on what to keep in mind technical interview. I am really excited if I am given an opportunity to work for </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>google</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> and hopefully towards an area where bitmap indices, compression schemes, or just database stuff is used, since I have some experience with that. But I have to survive the interview process first ;).</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="NoSpacing"/></w:pPr></w:p><w:p><w:pPr><w:pStyle w:val="NoSpacing"/></w:pPr><w:r><w:t>All the best,</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="NoSpacing"/><w:rPr><w:sz w:val="17"/><w:szCs w:val="17"/></w:rPr></w:pPr><w:r><w:t>Eduardo Gutarra</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$d.Range(0, 0).InsertXML($bodyXml)

Write-Output $d.Paragraphs.Count
